# Clear the E3:E22 "controle de integridade do banco" column on the
# "Avaliação" sheet (Grupo 3 / E column progress marks), matching the
# commit that removed the interim "m"/"r" markers and numeric marks.
$wb = $excel.ActiveWorkbook

$wsAval = $wb.Worksheets.Item("Avaliação")
$wsAcomp = $wb.Worksheets.Item("Acompanhamento")

$wsAval.Range("E3:E22").ClearContents()

# Update selection / view state to match the authored diff.
$wsAval.Application.ActiveWindow.ScrollRow = 1
$wsAval.Range("C20").Select()

$wsAcomp.Application.ActiveWindow.ScrollRow = 1
